$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 470
$ws.Range("F4").Value = 215
$ws.Range("F5").Value = 74
$ws.Range("F6").Value = 11
$ws.Range("F10").Value = 1311
$ws.Range("F12").Value = 1093
$ws.Range("F13").Value = 24
$ws.Range("F16").Value = 106
$ws.Range("F17").Value = 246
$ws.Range("F18").Value = 1665
$ws.Range("F21").Value = 228
$ws.Range("F22").Value = 2388
$ws.Range("F23").Value = 10
$ws.Range("F24").Value = 404
$ws.Range("F26").Value = 932
$ws.Range("F27").Value = 1215
$ws.Range("F30").Value = 2827
$ws.Range("F31").Value = 1624
$ws.Range("F32").Value = 84
$ws.Range("F33").Value = 118
$ws.Range("F34").Value = 672
$ws.Range("F35").Value = 866
$ws.Range("F36").Value = 1843
$ws.Range("F37").Value = 894
$ws.Range("F38").Value = 1861
$ws.Range("F41").Value = 844
$ws.Range("F42").Value = 41
$ws.Range("F43").Value = 874
$ws.Range("F44").Value = 799
$ws.Range("F45").Value = 1023
$ws.Range("F46").Value = 95
$ws.Range("F47").Value = 439
$ws.Range("F48").Value = 221
$ws.Range("F49").Value = 3347
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 11
$ws.Range("F12").Value = 803
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 470
$ws.Range("F4").Value = 215
$ws.Range("F5").Value = 74
$ws.Range("F6").Value = 11
$ws.Range("F8").Value = 11
$ws.Range("F12").Value = 1311
$ws.Range("F14").Value = 1093
$ws.Range("F15").Value = 24
$ws.Range("F18").Value = 106
$ws.Range("F19").Value = 246
$ws.Range("F20").Value = 1665
$ws.Range("F23").Value = 228
$ws.Range("F24").Value = 2389
$ws.Range("F25").Value = 404
$ws.Range("F27").Value = 1215
$ws.Range("F28").Value = 2827
$ws.Range("F29").Value = 1624
$ws.Range("F30").Value = 84
$ws.Range("F31").Value = 118
$ws.Range("F32").Value = 803
$ws.Range("F34").Value = 672
$ws.Range("F35").Value = 866
$ws.Range("F36").Value = 1843
$ws.Range("F38").Value = 894
$ws.Range("F39").Value = 1861
$ws.Range("F40").Value = 844
$ws.Range("F41").Value = 874
$ws.Range("F42").Value = 799
$ws.Range("F43").Value = 1023
$ws.Range("F44").Value = 95
$ws.Range("F45").Value = 439
$ws.Range("F47").Value = 221
$ws.Range("F48").Value = 3347
